$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number need an explicit
# text format first, otherwise Excel auto-converts the typed string
# into a real number (the source data keeps these as text cells,
# e.g. "0.571", "528.98", matching the original column formatting).
$ws.Range("D2").Value = "59.350.28"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").Value = "2.641.93"
$ws.Range("E3").Value = "  -0.23%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "528.98"
$ws.Range("E5").Value = "  +1.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.27"
$ws.Range("E6").Value = "  -0.70%  "
$ws.Range("E7").Value = "  -0.24%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.571"
$ws.Range("E8").Value = "  +0.21%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.67"
$ws.Range("E9").Value = "  -3.53%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.105"
$ws.Range("E10").Value = "  +2.38%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.339"
$ws.Range("E11").Value = "  +1.22%  "
$ws.Range("E12").Value = "  +1.11%  "
$ws.Range("D13").Value = "3.110.87"
$ws.Range("E13").Value = "  -0.22%  "
$ws.Range("D14").Value = "59.286.38"
$ws.Range("E14").Value = "  +0.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.04"
$ws.Range("E15").Value = "  +0.24%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000137"
$ws.Range("E16").Value = "  +1.24%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.672.22"
$ws.Range("E17").Value = "  +0.95%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "342.27"
$ws.Range("E18").Value = "  +0.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.47"
$ws.Range("E19").Value = "  +1.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.63"
$ws.Range("E20").Value = "  +3.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.37"
$ws.Range("E21").Value = "  +1.53%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  +0.35%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.68"
$ws.Range("E23").Value = "  +3.83%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.420"
$ws.Range("E24").Value = "  +2.27%  "
$ws.Range("E25").Value = "  +0.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.997"
$ws.Range("E26").Value = "  -0.46%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.28"
$ws.Range("E27").Value = "  +2.13%  "
$ws.Range("D28").Value = "0.0₃0803"
$ws.Range("E28").Value = "  +0.43%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.48"
$ws.Range("E29").Value = "  -2.67%  "
$ws.Range("E30").Value = "  -0.08%  "
$ws.Range("E31").Value = "  +2.29%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.03"
$ws.Range("E32").Value = "  +1.67%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "150.15"
$ws.Range("E33").Value = "  +0.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.23"
$ws.Range("E34").Value = "  +1.97%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.21"
$ws.Range("E35").Value = "  +1.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.899"
$ws.Range("E36").Value = "  +0.15%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.871"
$ws.Range("E37").Value = "  -0.88%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.50"
$ws.Range("E38").Value = "  +1.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.64"
$ws.Range("E39").Value = "  -0.78%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.67"
$ws.Range("E40").Value = "  +2.68%  "
$ws.Range("E41").Value = "  -0.44%  "
$ws.Range("E42").Value = "  +0.35%  "
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "272.98"
$ws.Range("E43").Value = "  -0.78%  "
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.603"
$ws.Range("E44").Value = "  -4.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.46"
$ws.Range("E45").Value = "  -1.47%  "
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0540"
$ws.Range("E46").Value = "  +0.94%  "
$ws.Range("B47").Value = "WhiteBITCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.66"
$ws.Range("E47").Value = "  +1.43%  "
$ws.Range("D48").Value = "2.042.87"
$ws.Range("E48").Value = "  -0.68%  "
$ws.Range("E49").Value = "  +0.16%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0231"
$ws.Range("E50").Value = "  +0.97%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.03"
$ws.Range("E51").Value = "  +0.83%  "
